$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price (column D) value still looks like a plain number need to be
# forced to Text format first, otherwise Excel auto-converts the literal string into a
# numeric cell (the source workbook stores these prices as text, e.g. "0.999", "568.54").
$textPriceCells = @("D4", "D5", "D6", "D7", "D10", "D11", "D13", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D29", "D30", "D32", "D33", "D35", "D36", "D39", "D41", "D42", "D43", "D46", "D49", "D50")
foreach ($ref in $textPriceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.017.27'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '2.541.69'
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '569.04'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").Value = '147.50'
$ws.Range("E6").Value = '  +3.83%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").Value = '2.540.38'
$ws.Range("E9").Value = '  +4.41%  '
$ws.Range("D10").Value = '0.105'
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").Value = '5.62'
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("E14").Value = '  +3.18%  '
$ws.Range("D15").Value = '2.994.96'
$ws.Range("E15").Value = '  +4.49%  '
$ws.Range("D16").Value = '62.938.76'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '2.545.51'
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("D19").Value = '11.49'
$ws.Range("E19").Value = '  +2.57%  '
$ws.Range("D20").Value = '333.46'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").Value = '4.29'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '6.73'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '64.81'
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  -2.69%  '
$ws.Range("E26").Value = '  +5.32%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +11.39%  '
$ws.Range("D29").Value = '8.39'
$ws.Range("E29").Value = '  +3.15%  '
$ws.Range("D30").Value = '7.17'
$ws.Range("E30").Value = '  +7.41%  '
$ws.Range("D31").Value = '0.0₃0814'
$ws.Range("E31").Value = '  +3.03%  '
$ws.Range("D32").Value = '1.85'
$ws.Range("E32").Value = '  +1.65%  '
$ws.Range("D33").Value = '177.34'
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("E34").Value = '  +6.44%  '
$ws.Range("D35").Value = '414.07'
$ws.Range("E35").Value = '  +11.70%  '
$ws.Range("D36").Value = '0.397'
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '4.36'
$ws.Range("E39").Value = '  -2.69%  '
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").Value = '39.17'
$ws.Range("E42").Value = '  -3.01%  '
$ws.Range("D43").Value = '151.35'
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").Value = '0.604'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D49").Value = '0.0235'
$ws.Range("E49").Value = '  +4.30%  '
$ws.Range("D50").Value = '18.31'
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("E51").Value = '  +3.23%  '

Write-Output "Updated cryptos list"
